$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.355.00"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "1.881.00"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'0.7115"
$ws.Range("E5").Value = "  +0.00%  "
$ws.Range("D6").Value = "'242.69"
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").Value = "'0.08024"
$ws.Range("E8").Value = "  +3.01%  "
$ws.Range("D9").Value = "'0.3139"
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("D10").Value = "'25.09"
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("D11").Value = "'0.08330"
$ws.Range("E11").Value = "  -1.43%  "
$ws.Range("D12").Value = "1.894.40"
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("D13").Value = "'5.262"
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("D14").Value = "'94.84"
$ws.Range("E14").Value = "  +4.05%  "
$ws.Range("D15").Value = "'0.7184"
$ws.Range("E15").Value = "  +0.93%  "
$ws.Range("D16").Value = "'6.373"
$ws.Range("E16").Value = "  +5.39%  "
$ws.Range("D17").Value = "'0.000008664"
$ws.Range("E17").Value = "  +5.38%  "
$ws.Range("D18").Value = "29.367.25"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "'243.33"
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").Value = "2.154.80"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").Value = "'13.34"
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").Value = "'1.002"
$ws.Range("E22").Value = "  +0.17%  "
$ws.Range("D23").Value = "'7.836"
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").Value = "'0.1575"
$ws.Range("E25").Value = "  -1.90%  "
$ws.Range("D26").Value = "'163.47"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "'9.093"
$ws.Range("E27").Value = "  +0.52%  "
$ws.Range("D28").Value = "'18.60"
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("D29").Value = "'1.511"
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "'4.433"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("D31").Value = "'4.361"
$ws.Range("E31").Value = "  +1.40%  "
$ws.Range("E32").Value = "  -6.50%  "
$ws.Range("D33").Value = "'0.05389"
$ws.Range("E33").Value = "  +2.15%  "
$ws.Range("D34").Value = "'1.943"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "'0.7771"
$ws.Range("E35").Value = "  +4.29%  "
$ws.Range("D36").Value = "'1.180"
$ws.Range("E36").Value = "  +0.20%  "
$ws.Range("E37").Value = "  -0.30%  "
$ws.Range("D38").Value = "'0.01886"
$ws.Range("E38").Value = "  +0.93%  "
$ws.Range("D39").Value = "1.270.81"
$ws.Range("E39").Value = "  +4.56%  "
$ws.Range("D40").Value = "'2.745"
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("D41").Value = "'6.527"
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("D42").Value = "'0.9199"
$ws.Range("E42").Value = "  +3.48%  "
$ws.Range("D43").Value = "'113.54"
$ws.Range("E43").Value = "  +4.24%  "
$ws.Range("D44").Value = "'74.58"
$ws.Range("E44").Value = "  +2.55%  "
$ws.Range("D45").Value = "'1.001"
$ws.Range("E45").Value = "  +0.09%  "
$ws.Range("D46").Value = "2.044.56"
$ws.Range("E46").Value = "  +1.34%  "
$ws.Range("E47").Value = "  +4.49%  "
$ws.Range("D48").Value = "'1.811"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("D49").Value = "'0.5224"
$ws.Range("E49").Value = "  +0.29%  "
$ws.Range("D50").Value = "'9.578"
$ws.Range("E50").Value = "  +2.19%  "
$ws.Range("D51").Value = "'0.4381"
$ws.Range("E51").Value = "  +1.45%  "
